$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (bold font, border, centered/top alignment) from A215
# down across the new A216:A270 range so it matches existing "index" column styling.
$ws.Range("A215").Copy()
$ws.Range("A216:A270").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(216, 1).Value = 214
$ws.Cells.Item(216, 2).Value = 215
$ws.Cells.Item(216, 3).Value = 1
$ws.Cells.Item(216, 4).Value = 'Начал взаимодействие с консультантом!'
$ws.Cells.Item(216, 5).Value = '11/06/2023 17:45:07'

$ws.Cells.Item(217, 1).Value = 215
$ws.Cells.Item(217, 2).Value = 216
$ws.Cells.Item(217, 3).Value = 1
$ws.Cells.Item(217, 4).Value = 'Рассказал партнеру о важности тестирования'
$ws.Cells.Item(217, 5).Value = '11/06/2023 17:45:45'

$ws.Cells.Item(218, 1).Value = 216
$ws.Cells.Item(218, 2).Value = 217
$ws.Cells.Item(218, 3).Value = 2
$ws.Cells.Item(218, 4).Value = 'Начал взаимодействие с консультантом!'
$ws.Cells.Item(218, 5).Value = '11/06/2023 17:56:53'

$ws.Cells.Item(219, 1).Value = 217
$ws.Cells.Item(219, 2).Value = 218
$ws.Cells.Item(219, 3).Value = 3
$ws.Cells.Item(219, 4).Value = 'Начал взаимодействие с консультантом!'
$ws.Cells.Item(219, 5).Value = '12/06/2023 06:34:59'

$ws.Cells.Item(220, 1).Value = 218
$ws.Cells.Item(220, 2).Value = 219
$ws.Cells.Item(220, 3).Value = 3
$ws.Cells.Item(220, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(220, 5).Value = '12/06/2023 06:49:07'

$ws.Cells.Item(221, 1).Value = 219
$ws.Cells.Item(221, 2).Value = 220
$ws.Cells.Item(221, 3).Value = 1
$ws.Cells.Item(221, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(221, 5).Value = '12/06/2023 06:49:12'

$ws.Cells.Item(222, 1).Value = 220
$ws.Cells.Item(222, 2).Value = 221
$ws.Cells.Item(222, 3).Value = 3
$ws.Cells.Item(222, 4).Value = 'Завершил тест sogi_assessment!'
$ws.Cells.Item(222, 5).Value = '12/06/2023 06:49:29'

$ws.Cells.Item(223, 1).Value = 221
$ws.Cells.Item(223, 2).Value = 222
$ws.Cells.Item(223, 3).Value = 1
$ws.Cells.Item(223, 4).Value = 'Завершил тест sogi_assessment!'
$ws.Cells.Item(223, 5).Value = '12/06/2023 06:49:36'

$ws.Cells.Item(224, 1).Value = 222
$ws.Cells.Item(224, 2).Value = 223
$ws.Cells.Item(224, 3).Value = 3
$ws.Cells.Item(224, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(224, 5).Value = '12/06/2023 06:49:43'

$ws.Cells.Item(225, 1).Value = 223
$ws.Cells.Item(225, 2).Value = 224
$ws.Cells.Item(225, 3).Value = 1
$ws.Cells.Item(225, 4).Value = 'Начал тест: hiv_risk_assessment'
$ws.Cells.Item(225, 5).Value = '12/06/2023 06:49:45'

$ws.Cells.Item(226, 1).Value = 224
$ws.Cells.Item(226, 2).Value = 225
$ws.Cells.Item(226, 3).Value = 3
$ws.Cells.Item(226, 4).Value = 'Завершил тест pkp_assessment!'
$ws.Cells.Item(226, 5).Value = '12/06/2023 06:49:52'

$ws.Cells.Item(227, 1).Value = 225
$ws.Cells.Item(227, 2).Value = 226
$ws.Cells.Item(227, 3).Value = 1
$ws.Cells.Item(227, 4).Value = 'Завершил тест hiv_risk_assessment!'
$ws.Cells.Item(227, 5).Value = '12/06/2023 06:50:22'

$ws.Cells.Item(228, 1).Value = 226
$ws.Cells.Item(228, 2).Value = 227
$ws.Cells.Item(228, 3).Value = 1
$ws.Cells.Item(228, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(228, 5).Value = '12/06/2023 06:50:29'

$ws.Cells.Item(229, 1).Value = 227
$ws.Cells.Item(229, 2).Value = 228
$ws.Cells.Item(229, 3).Value = 3
$ws.Cells.Item(229, 4).Value = 'Начал тест: understanding_PLHIV_assessment'
$ws.Cells.Item(229, 5).Value = '12/06/2023 06:50:40'

$ws.Cells.Item(230, 1).Value = 228
$ws.Cells.Item(230, 2).Value = 229
$ws.Cells.Item(230, 3).Value = 1
$ws.Cells.Item(230, 4).Value = 'Завершил тест pkp_assessment!'
$ws.Cells.Item(230, 5).Value = '12/06/2023 06:50:51'

$ws.Cells.Item(231, 1).Value = 229
$ws.Cells.Item(231, 2).Value = 230
$ws.Cells.Item(231, 3).Value = 1
$ws.Cells.Item(231, 4).Value = 'Начал тест: understanding_PLHIV_assessment'
$ws.Cells.Item(231, 5).Value = '12/06/2023 06:51:14'

$ws.Cells.Item(232, 1).Value = 230
$ws.Cells.Item(232, 2).Value = 231
$ws.Cells.Item(232, 3).Value = 3
$ws.Cells.Item(232, 4).Value = 'Завершил тест understanding_PLHIV_assessment!'
$ws.Cells.Item(232, 5).Value = '12/06/2023 06:51:59'

$ws.Cells.Item(233, 1).Value = 231
$ws.Cells.Item(233, 2).Value = 232
$ws.Cells.Item(233, 3).Value = 3
$ws.Cells.Item(233, 4).Value = 'Начал тест: hiv_risk_assessment'
$ws.Cells.Item(233, 5).Value = '12/06/2023 06:52:14'

$ws.Cells.Item(234, 1).Value = 232
$ws.Cells.Item(234, 2).Value = 233
$ws.Cells.Item(234, 3).Value = 3
$ws.Cells.Item(234, 4).Value = 'Завершил тест hiv_risk_assessment!'
$ws.Cells.Item(234, 5).Value = '12/06/2023 06:52:44'

$ws.Cells.Item(235, 1).Value = 233
$ws.Cells.Item(235, 2).Value = 234
$ws.Cells.Item(235, 3).Value = 6
$ws.Cells.Item(235, 4).Value = 'Успешно добавлен в базу!'
$ws.Cells.Item(235, 5).Value = '12/06/2023 07:07:13'

$ws.Cells.Item(236, 1).Value = 234
$ws.Cells.Item(236, 2).Value = 235
$ws.Cells.Item(236, 3).Value = 1
$ws.Cells.Item(236, 4).Value = 'Начал взаимодействие с консультантом!'
$ws.Cells.Item(236, 5).Value = '12/06/2023 07:20:13'

$ws.Cells.Item(237, 1).Value = 235
$ws.Cells.Item(237, 2).Value = 236
$ws.Cells.Item(237, 3).Value = 1
$ws.Cells.Item(237, 4).Value = 'Попытался заказать тест на ВИЧ!'
$ws.Cells.Item(237, 5).Value = '12/06/2023 07:20:22'

$ws.Cells.Item(238, 1).Value = 236
$ws.Cells.Item(238, 2).Value = 237
$ws.Cells.Item(238, 3).Value = 1
$ws.Cells.Item(238, 4).Value = 'Начал тест: hiv_knowledge_assessment'
$ws.Cells.Item(238, 5).Value = '12/06/2023 07:20:42'

$ws.Cells.Item(239, 1).Value = 237
$ws.Cells.Item(239, 2).Value = 238
$ws.Cells.Item(239, 3).Value = 1
$ws.Cells.Item(239, 4).Value = 'Попытался заказать тест на ВИЧ!'
$ws.Cells.Item(239, 5).Value = '12/06/2023 07:20:58'

$ws.Cells.Item(240, 1).Value = 238
$ws.Cells.Item(240, 2).Value = 239
$ws.Cells.Item(240, 3).Value = 1
$ws.Cells.Item(240, 4).Value = 'Начал взаимодействие с консультантом!'
$ws.Cells.Item(240, 5).Value = '12/06/2023 07:21:15'

$ws.Cells.Item(241, 1).Value = 239
$ws.Cells.Item(241, 2).Value = 240
$ws.Cells.Item(241, 3).Value = 6
$ws.Cells.Item(241, 4).Value = 'Начал взаимодействие с консультантом!'
$ws.Cells.Item(241, 5).Value = '12/06/2023 07:21:58'

$ws.Cells.Item(242, 1).Value = 240
$ws.Cells.Item(242, 2).Value = 241
$ws.Cells.Item(242, 3).Value = 1
$ws.Cells.Item(242, 4).Value = 'Успешно добавлен в базу!'
$ws.Cells.Item(242, 5).Value = '12/06/2023 07:42:36'

$ws.Cells.Item(243, 1).Value = 241
$ws.Cells.Item(243, 2).Value = 242
$ws.Cells.Item(243, 3).Value = 1
$ws.Cells.Item(243, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(243, 5).Value = '16/06/2023 14:42:10'

$ws.Cells.Item(244, 1).Value = 242
$ws.Cells.Item(244, 2).Value = 243
$ws.Cells.Item(244, 3).Value = 1
$ws.Cells.Item(244, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(244, 5).Value = '16/06/2023 14:44:57'

$ws.Cells.Item(245, 1).Value = 243
$ws.Cells.Item(245, 2).Value = 244
$ws.Cells.Item(245, 3).Value = 1
$ws.Cells.Item(245, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(245, 5).Value = '16/06/2023 15:20:24'

$ws.Cells.Item(246, 1).Value = 244
$ws.Cells.Item(246, 2).Value = 245
$ws.Cells.Item(246, 3).Value = 1
$ws.Cells.Item(246, 4).Value = 'Завершил тест pkp_assessment!'
$ws.Cells.Item(246, 5).Value = '16/06/2023 15:24:16'

$ws.Cells.Item(247, 1).Value = 245
$ws.Cells.Item(247, 2).Value = 246
$ws.Cells.Item(247, 3).Value = 1
$ws.Cells.Item(247, 4).Value = 'Начал тест: hiv_knowledge_assessment'
$ws.Cells.Item(247, 5).Value = '16/06/2023 15:26:20'

$ws.Cells.Item(248, 1).Value = 246
$ws.Cells.Item(248, 2).Value = 247
$ws.Cells.Item(248, 3).Value = 1
$ws.Cells.Item(248, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(248, 5).Value = '16/06/2023 15:27:53'

$ws.Cells.Item(249, 1).Value = 247
$ws.Cells.Item(249, 2).Value = 248
$ws.Cells.Item(249, 3).Value = 1
$ws.Cells.Item(249, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(249, 5).Value = '16/06/2023 15:31:10'

$ws.Cells.Item(250, 1).Value = 248
$ws.Cells.Item(250, 2).Value = 249
$ws.Cells.Item(250, 3).Value = 1
$ws.Cells.Item(250, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(250, 5).Value = '16/06/2023 15:32:09'

$ws.Cells.Item(251, 1).Value = 249
$ws.Cells.Item(251, 2).Value = 250
$ws.Cells.Item(251, 3).Value = 1
$ws.Cells.Item(251, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(251, 5).Value = '16/06/2023 15:36:07'

$ws.Cells.Item(252, 1).Value = 250
$ws.Cells.Item(252, 2).Value = 251
$ws.Cells.Item(252, 3).Value = 1
$ws.Cells.Item(252, 4).Value = 'Начал тест: hiv_risk_assessment'
$ws.Cells.Item(252, 5).Value = '16/06/2023 15:36:46'

$ws.Cells.Item(253, 1).Value = 251
$ws.Cells.Item(253, 2).Value = 252
$ws.Cells.Item(253, 3).Value = 1
$ws.Cells.Item(253, 4).Value = 'Начал тест: hiv_risk_assessment'
$ws.Cells.Item(253, 5).Value = '16/06/2023 16:07:15'

$ws.Cells.Item(254, 1).Value = 252
$ws.Cells.Item(254, 2).Value = 253
$ws.Cells.Item(254, 3).Value = 1
$ws.Cells.Item(254, 4).Value = 'Завершил тест hiv_risk_assessment!'
$ws.Cells.Item(254, 5).Value = '16/06/2023 16:07:58'

$ws.Cells.Item(255, 1).Value = 253
$ws.Cells.Item(255, 2).Value = 254
$ws.Cells.Item(255, 3).Value = 1
$ws.Cells.Item(255, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(255, 5).Value = '16/06/2023 16:08:13'

$ws.Cells.Item(256, 1).Value = 254
$ws.Cells.Item(256, 2).Value = 255
$ws.Cells.Item(256, 3).Value = 1
$ws.Cells.Item(256, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(256, 5).Value = '16/06/2023 16:09:40'

$ws.Cells.Item(257, 1).Value = 255
$ws.Cells.Item(257, 2).Value = 256
$ws.Cells.Item(257, 3).Value = 1
$ws.Cells.Item(257, 4).Value = 'Завершил тест pkp_assessment!'
$ws.Cells.Item(257, 5).Value = '16/06/2023 16:10:25'

$ws.Cells.Item(258, 1).Value = 256
$ws.Cells.Item(258, 2).Value = 257
$ws.Cells.Item(258, 3).Value = 1
$ws.Cells.Item(258, 4).Value = 'Начал тест: pkp_assessment'
$ws.Cells.Item(258, 5).Value = '16/06/2023 16:11:12'

$ws.Cells.Item(259, 1).Value = 257
$ws.Cells.Item(259, 2).Value = 258
$ws.Cells.Item(259, 3).Value = 1
$ws.Cells.Item(259, 4).Value = 'Завершил тест pkp_assessment!'
$ws.Cells.Item(259, 5).Value = '16/06/2023 16:12:07'

$ws.Cells.Item(260, 1).Value = 258
$ws.Cells.Item(260, 2).Value = 259
$ws.Cells.Item(260, 3).Value = 1
$ws.Cells.Item(260, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(260, 5).Value = '16/06/2023 16:29:18'

$ws.Cells.Item(261, 1).Value = 259
$ws.Cells.Item(261, 2).Value = 260
$ws.Cells.Item(261, 3).Value = 1
$ws.Cells.Item(261, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(261, 5).Value = '16/06/2023 16:30:47'

$ws.Cells.Item(262, 1).Value = 260
$ws.Cells.Item(262, 2).Value = 261
$ws.Cells.Item(262, 3).Value = 1
$ws.Cells.Item(262, 4).Value = 'Завершил тест sogi_assessment!'
$ws.Cells.Item(262, 5).Value = '16/06/2023 16:31:41'

$ws.Cells.Item(263, 1).Value = 261
$ws.Cells.Item(263, 2).Value = 262
$ws.Cells.Item(263, 3).Value = 1
$ws.Cells.Item(263, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(263, 5).Value = '16/06/2023 17:57:05'

$ws.Cells.Item(264, 1).Value = 262
$ws.Cells.Item(264, 2).Value = 263
$ws.Cells.Item(264, 3).Value = 1
$ws.Cells.Item(264, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(264, 5).Value = '16/06/2023 17:59:27'

$ws.Cells.Item(265, 1).Value = 263
$ws.Cells.Item(265, 2).Value = 264
$ws.Cells.Item(265, 3).Value = 1
$ws.Cells.Item(265, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(265, 5).Value = '16/06/2023 18:00:41'

$ws.Cells.Item(266, 1).Value = 264
$ws.Cells.Item(266, 2).Value = 265
$ws.Cells.Item(266, 3).Value = 1
$ws.Cells.Item(266, 4).Value = 'Завершил тест sogi_assessment!'
$ws.Cells.Item(266, 5).Value = '16/06/2023 18:02:52'

$ws.Cells.Item(267, 1).Value = 265
$ws.Cells.Item(267, 2).Value = 266
$ws.Cells.Item(267, 3).Value = 1
$ws.Cells.Item(267, 4).Value = 'Начал тест: sogi_assessment'
$ws.Cells.Item(267, 5).Value = '16/06/2023 18:03:31'

$ws.Cells.Item(268, 1).Value = 266
$ws.Cells.Item(268, 2).Value = 267
$ws.Cells.Item(268, 3).Value = 1
$ws.Cells.Item(268, 4).Value = 'Завершил тест sogi_assessment!'
$ws.Cells.Item(268, 5).Value = '16/06/2023 18:04:25'

$ws.Cells.Item(269, 1).Value = 267
$ws.Cells.Item(269, 2).Value = 268
$ws.Cells.Item(269, 3).Value = 1
$ws.Cells.Item(269, 4).Value = 'Начал тест: hiv_risk_assessment'
$ws.Cells.Item(269, 5).Value = '16/06/2023 18:04:44'

$ws.Cells.Item(270, 1).Value = 268
$ws.Cells.Item(270, 2).Value = 269
$ws.Cells.Item(270, 3).Value = 1
$ws.Cells.Item(270, 4).Value = 'Завершил тест hiv_risk_assessment!'
$ws.Cells.Item(270, 5).Value = '16/06/2023 18:05:25'
